{"js": "const body = context.document.body;\n\n// 1) \"Versi\" + \"on\" -> merge into a single run \"Version\" (text unchanged,\n//    but the two runs collapse into one once rewritten).\nlet results = body.search(\"Version\", { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\nresults.items[0].insertText(\"Version\", Word.InsertLocation.replace);\nawait context.sync();\n\n// 2) \" 2\" -> \" 1.\" (the run that used to hold just the version number now\n//    also carries the trailing period).\nresults = body.search(\" 2\", { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\nresults.items[0].insertText(\" 1.\", Word.InsertLocation.replace);\nawait context.sync();\n\n// 3) Remove the now-redundant trailing \".\" run that sat after the bookmark.\nresults = body.search(\".\", { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\nconst trailingPeriod = results.items[results.items.length - 1];\ntrailingPeriod.insertText(\"\", Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# 1) \"Versi\" + \"on\" -> merge into a single run \"Version\" (the two original\n#    runs collapse into one once Word rewrites the matched text).\n$find1 = $d.Content.Find\n$find1.Execute(\"Version\", $false, $false, $false, $false, $false, $true, 1, $false, \"Version\", 2) | Out-Null\n\n# 2) \" 2\" -> \" 1.\" (this run now also carries the trailing period).\n$find2 = $d.Content.Find\n$find2.Execute(\" 2\", $false, $false, $false, $false, $false, $true, 1, $false, \" 1.\", 2) | Out-Null\n\n# 3) Drop the now-redundant \".\" that used to sit in its own run right after\n#    the bookmark, at the very end of the paragraph.\n$para = $d.Paragraphs(1)\n$pr = $para.Range\n$lastCharIndex = $pr.Characters.Count - 1\n$pr.Characters($lastCharIndex).Delete()\n"}
